$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 528; this shifts the existing rows 528:632 down to 529:633
# (and updates the sheet dimension to A1:R633 automatically).
$ws.Rows(528).Insert()

# Populate the newly inserted row 528 with the new record.
$ws.Range("A528").Value = 4
$ws.Range("B528").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C528").Value = "Los Lagos"
$ws.Range("D528").Value = 45258
$ws.Range("D528").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E528").Value = 10
$ws.Range("F528").Value = 100112023
$ws.Range("G528").Value = "Brócoli"
$ws.Range("H528").Value = "Sin especificar"
$ws.Range("I528").Value = "Primera"
$ws.Range("J528").Value = 1500
$ws.Range("K528").Value = 1600
$ws.Range("L528").Value = 1600
$ws.Range("M528").Value = 1600
$ws.Range("N528").Value = "$/unidad"
$ws.Range("O528").Value = "Región Metropolitana"
$ws.Range("P528").Value = 1600
$ws.Range("Q528").Value = 1
$ws.Range("R528").Value = "Hortaliza"
